# Apply updated KNN imputation results to Sheet1.
# This mirrors an upstream "Update Name of Algo" commit where the
# underlying result_data_KNN.xlsx values were refreshed for several cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -10.695
$ws.Range("D3").Value = -7.264999999999999
$ws.Range("E8").Value = 16.778
$ws.Range("E11").Value = 17.103
$ws.Range("A12").Value = -21.683
$ws.Range("C14").Value = -12.465
$ws.Range("E14").Value = 17.186
$ws.Range("E15").Value = 15.927
$ws.Range("C26").Value = -12.878
$ws.Range("D30").Value = -7.257
$ws.Range("C31").Value = -12.597
$ws.Range("A32").Value = -21.79799999999999
$ws.Range("C35").Value = -12.762
$ws.Range("A36").Value = -20.178
$ws.Range("E36").Value = 16.44
$ws.Range("C37").Value = -13.578
$ws.Range("A38").Value = -19.741
$ws.Range("D44").Value = -7.746
$ws.Range("C45").Value = -12.883
$ws.Range("A46").Value = -21.924
$ws.Range("A54").Value = -22.155
$ws.Range("A55").Value = -22.278
$ws.Range("C57").Value = -13.829
$ws.Range("D58").Value = -8.040000000000001
$ws.Range("E64").Value = 17.185
$ws.Range("A67").Value = -21.493
$ws.Range("A69").Value = -21.544
$ws.Range("A72").Value = -21.445
$ws.Range("D84").Value = -8.172000000000001
$ws.Range("D89").Value = -7.233999999999999
$ws.Range("E89").Value = 17.078
$ws.Range("A91").Value = -21.587
$ws.Range("D91").Value = -6.910000000000001
$ws.Range("D92").Value = -6.737
$ws.Range("A99").Value = -20.428
$ws.Range("C100").Value = -12.156
$ws.Range("C102").Value = -13.774
$ws.Range("D102").Value = -7.764

$wb.Save()
